# Add five new rows (3-7) of CvLAC data below the existing header/data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dept = "Departamento de Ingeniería de Sistemas e Industrial"

# --- Names / department for the five new people (written row by row) ---
$ws.Range("A3").Value = 123123123
$ws.Range("B3").Value = "Wilson"
$ws.Range("C3").Value = "Adarme Jaimes"
$ws.Range("D3").Value = $dept

$ws.Range("A4").Value = 123123123
$ws.Range("B4").Value = "Jairo Hernán"
$ws.Range("C4").Value = "Aponte Melo"
$ws.Range("D4").Value = $dept

$ws.Range("A5").Value = 123123123
$ws.Range("B5").Value = "Emiliano"
$ws.Range("C5").Value = "Barreto Hernández"
$ws.Range("D5").Value = $dept

$ws.Range("A6").Value = 123123123
$ws.Range("B6").Value = "Libia Denise"
$ws.Range("C6").Value = "Cangrejo Aljure"
$ws.Range("D6").Value = $dept

$ws.Range("A7").Value = 123123123
$ws.Range("B7").Value = "Ismael"
$ws.Range("C7").Value = "Castañeda Fuentes"
$ws.Range("D7").Value = $dept

# --- CvLAC links (plain text, not hyperlinks, for rows 3,4,6,7) ---
$ws.Range("E3").Value = "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0000439185"
$ws.Range("E4").Value = "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0001333865"
$ws.Range("E6").Value = "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0001370358"
$ws.Range("E7").Value = "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0000199087"

# Row 5's link is an actual hyperlink (mirrors the one already on E2).
$ws.Hyperlinks.Add($ws.Range("E5"), "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0000025410")

# Move the selection the way the author left it.
[void]$ws.Range("A5").Select()
